$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the contents of columns B and C for each data row (1..32).
# The college-name and code columns were mismapped; this fixes the mapping.
for ($r = 1; $r -le 32; $r++) {
    $bCell = $ws.Cells.Item($r, 2)
    $cCell = $ws.Cells.Item($r, 3)
    $bVal = $bCell.Value()
    $cVal = $cCell.Value()
    $bCell.Value = $cVal
    $cCell.Value = $bVal
}

# Columns B and C also swap their widths/best-fit sizing to match the new content.
$ws.Columns(2).ColumnWidth = 16.57
$ws.Columns(3).ColumnWidth = 26.43

# Update the active selection as recorded in the saved workbook.
$ws.Range("I8").Select()
